# "Add stage 3 for B1-B2"
# Adds a third pair of TxHash values (stage 3) to the B1 and B2 sheets,
# and updates the active sheet/selection state to reflect where the
# author was working (B2 sheet, cell G11) instead of the Info sheet.

$wb = $excel.ActiveWorkbook

# --- B1 sheet: new stage-3 hashes in A2/A3 ---
$wsB1 = $wb.Worksheets.Item("B1")
$wsB1.Range("A2").Value = "7BAA7642C2250AE32DB955C777C19B2C470A8D839B74BB4566C7E4314B5CDF30"
$wsB1.Range("A3").Value = "0402D337EDA48FC520830DE8096AB6A429BCA96C1230FFE9A8D2AFD844BADF2A"

# --- B2 sheet: new stage-3 hashes in A2/A3 ---
$wsB2 = $wb.Worksheets.Item("B2")
$wsB2.Range("A2").Value = "A47AE2BDEC9179D0169FE5D8F5438D8FC6C1B9836791B5B55086CB127066F226"
$wsB2.Range("A3").Value = "5F3CE35BE89900A27210E7179665E0F07D65979E60ED59E39E776679987C454A"

# --- Update selections to match the author's last-saved cursor state ---
$wsInfo = $wb.Worksheets.Item("Info")
[void]$wsInfo.Activate()
[void]$wsInfo.Range("B2").Select()

[void]$wsB1.Range("A4").Select()

[void]$wsB2.Activate()
[void]$wsB2.Range("G11").Select()
